# Updating the Base_Map_Closed to the 2019.2 PreScan Version
# TAconstant.xlsx: the GoalTime constant (row 19, "BMW_X5_SUV_1/GoalTime")
# is updated from 10 to 1, and the sheet's selection ends on B20 (the cell
# below it) as it would after editing B19 and pressing Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the constant value in column B, row 19 (GoalTime): 10 -> 1
$ws.Range("B19").Value = 1

# Leave the selection on B20, matching the post-edit cursor position
$ws.Range("B20").Select()
